# Applies the "Updated cryptos list" data refresh to sheet1 (Coin/Link/Price/Volume table).
# Values are stored as text in the sheet (e.g. "27.027.30", "  -0.50%  "), so for any
# Price value that Excel would otherwise auto-parse as a number (dropping trailing
# zeros / switching to scientific notation, e.g. "0.0630" -> 6.3E-2), we temporarily
# force the cell to Text format, assign the literal string, then restore the cell
# style to Normal so the saved file keeps no explicit style on the cell (matching the
# original, unstyled data cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.049.34"
$ws.Range("E2").Value = "  -0.40%  "
# Row 3
$ws.Range("D3").Value = "1.623.27"
$ws.Range("E3").Value = "  -1.07%  "
# Row 4
$ws.Range("E4").Value = "  -0.13%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
# Row 7
$ws.Range("E7").Value = "  -0.13%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0630"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.70%  "
# Row 9
$ws.Range("E9").Value = "  -1.51%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
# Row 11
$ws.Range("E11").Value = "  +0.00%  "
# Row 12
$ws.Range("D12").Value = "1.850.51"
$ws.Range("E12").Value = "  -1.05%  "
# Row 13
$ws.Range("D13").Value = "1.624.03"
$ws.Range("E13").Value = "  -1.02%  "
# Row 14
$ws.Range("E14").Value = "  +0.33%  "
# Row 15
$ws.Range("E15").Value = "  -0.08%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.17%  "
# Row 17
$ws.Range("D17").Value = "27.020.72"
$ws.Range("E17").Value = "  -0.52%  "
# Row 18
$ws.Range("E18").Value = "  +0.66%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.48%  "
# Row 20
$ws.Range("E20").Value = "  -0.14%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.12%  "
# Row 22
$ws.Range("E22").Value = "  -0.95%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.36%  "
# Row 24
$ws.Range("E24").Value = "  -0.49%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
# Row 26
$ws.Range("E26").Value = "  -0.10%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "
# Row 28
$ws.Range("E28").Value = "  -2.63%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.59%  "
# Row 30
$ws.Range("E30").Value = "  +0.78%  "
# Row 31
$ws.Range("E31").Value = "  -0.73%  "
# Row 32
$ws.Range("E32").Value = "  -0.94%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.749"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +36.49%  "
# Row 34
$ws.Range("E34").Value = "  -0.07%  "
# Row 35
$ws.Range("D35").Value = "1.360.73"
$ws.Range("E35").Value = "  +4.33%  "
# Row 36
$ws.Range("E36").Value = "  +0.32%  "
# Row 37
$ws.Range("E37").Value = "  -0.81%  "
# Row 38
$ws.Range("E38").Value = "  +1.03%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
# Row 40
$ws.Range("E40").Value = "  -0.16%  "
# Row 41
$ws.Range("E41").Value = "  -1.27%  "
# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.85%  "
# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "
# Row 45
$ws.Range("D45").Value = "1.762.23"
$ws.Range("E45").Value = "  -1.03%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.887"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +32.50%  "
# Row 47
$ws.Range("E47").Value = "  -1.93%  "
# Row 48
$ws.Range("E48").Value = "  +2.84%  "
# Row 49
$ws.Range("E49").Value = "  -0.76%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.10%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0515"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.46%  "
